$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 7-8 (shifting old rows 7..17 down to 9..19)
$ws.Rows("7:8").Insert()

# --- Column A: new folder names ---
$ws.Range("A7").Formula = "codev-lc16-ecommerce-website"
$ws.Range("A8").Formula = "codev-lc18-vanilla-js"

# --- Columns B:F: identical static HTML template fragments used by every row ---
foreach ($col in @("B", "C", "D", "E", "F")) {
    $srcFormula = $ws.Range("$col`6").Formula
    $ws.Range("$col`7").Formula = $srcFormula
    $ws.Range("$col`8").Formula = $srcFormula
}

# --- Columns G:I: the shared CONCAT formulas ---
$ws.Range("G7").Formula = "=_xlfn.CONCAT(`$B7,`$A7,`$C7,`$A7,`$D7)"
$ws.Range("H7").Formula = "=_xlfn.CONCAT(`$A7,`$E7,`$A7,`$F7)"
$ws.Range("I7").Formula = "=_xlfn.CONCAT(`$G7,`$H7)"

$ws.Range("G8").Formula = "=_xlfn.CONCAT(`$B8,`$A8,`$C8,`$A8,`$D8)"
$ws.Range("H8").Formula = "=_xlfn.CONCAT(`$A8,`$E8,`$A8,`$F8)"
$ws.Range("I8").Formula = "=_xlfn.CONCAT(`$G8,`$H8)"

# --- Re-apply the row formatting the rest of the table uses ---
foreach ($r in 7, 8) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Borders.Item(7).LineStyle = 1
    $cellA.Borders.Item(10).LineStyle = 1

    foreach ($col in @("B", "C", "D", "E", "F")) {
        $ws.Range("$col$r").NumberFormat = "@"
    }

    $cellI = $ws.Cells.Item($r, 9)
    $cellI.Borders.Item(7).LineStyle = 1
    $cellI.Borders.Item(10).LineStyle = 1
    $cellI.Interior.Color = 65535
}

# --- Match the author's final selection state ---
$ws.Range("I7:I8").Select()
